# Auto-generated edit script applying numeric updates to Famfrit_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 326.2
$ws.Range("I33").Value = 325.92856
$ws.Range("K33").Value = 325.92856
$ws.Range("M33").Value = -96.92856

$ws.Range("H62").Value = 4416.5557
$ws.Range("I62").Value = 2681.9092
$ws.Range("K62").Value = 2681.9092
$ws.Range("M62").Value = -2057.9092

$ws.Range("H65").Value = 4416.5557
$ws.Range("I65").Value = 2681.9092
$ws.Range("K65").Value = 13409.546
$ws.Range("M65").Value = -10289.546

$ws.Range("H106").Value = 2530.5908
$ws.Range("I106").Value = 2341.2632
$ws.Range("K106").Value = 2341.2632
$ws.Range("M106").Value = -1710.2632

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1466.6666
$ws.Range("I5").Value = 1466.6666
$ws.Range("K5").Value = 1466.6666
$ws.Range("M5").Value = -1354.6666

$ws.Range("H32").Value = 10758486
$ws.Range("I32").Value = 13700904
$ws.Range("K32").Value = 13700904
$ws.Range("M32").Value = -13700617

$ws.Range("H45").Value = 1753.4667
$ws.Range("I45").Value = 1537.7
$ws.Range("K45").Value = 1537.7
$ws.Range("M45").Value = -1160.7

$ws.Range("H61").Value = 41670092
$ws.Range("I61").Value = 58824980
$ws.Range("J61").Value = 8216.143
$ws.Range("K61").Value = 58824980
$ws.Range("L61").Value = 8216.143
$ws.Range("M61").Value = -58824768
$ws.Range("N61").Value = -8640.143

$ws.Range("H63").Value = 4990.815
$ws.Range("I63").Value = 3237.85
$ws.Range("K63").Value = 3237.85
$ws.Range("M63").Value = -2551.85

$ws.Range("H66").Value = 4990.815
$ws.Range("I66").Value = 3237.85
$ws.Range("K66").Value = 16189.25
$ws.Range("M66").Value = -12757.25

$ws.Range("H88").Value = 10670.167
$ws.Range("I88").Value = 16449.715
$ws.Range("K88").Value = 16449.715
$ws.Range("M88").Value = -16043.715

$ws.Range("H91").Value = 10670.167
$ws.Range("I91").Value = 16449.715
$ws.Range("K91").Value = 16449.715
$ws.Range("M91").Value = -15045.715

$ws.Range("H97").Value = 1386.85
$ws.Range("I97").Value = 388.35715
$ws.Range("K97").Value = 388.35715
$ws.Range("M97").Value = 107.64285

$ws.Range("H102").Value = 4441.2
$ws.Range("I102").Value = 2100
$ws.Range("J102").Value = 6002
$ws.Range("K102").Value = 2100
$ws.Range("L102").Value = 6002
$ws.Range("M102").Value = -478
$ws.Range("N102").Value = -9246

$ws.Range("H110").Value = 15572.77
$ws.Range("I110").Value = 16682.479
$ws.Range("J110").Value = 7065
$ws.Range("K110").Value = 16682.479
$ws.Range("L110").Value = 7065
$ws.Range("M110").Value = -14637.479
$ws.Range("N110").Value = -11155

$ws.Range("H122").Value = 4670.0454
$ws.Range("I122").Value = 4148.5
$ws.Range("K122").Value = 12445.5
$ws.Range("M122").Value = -9995.5

$ws.Range("H132").Value = 52634644
$ws.Range("I132").Value = 3041.2354
$ws.Range("J132").Value = 500003260
$ws.Range("K132").Value = 9123.706200000001
$ws.Range("L132").Value = 1500009780
$ws.Range("M132").Value = -6593.706200000001
$ws.Range("N132").Value = -1500014840

$ws.Range("H136").Value = 41670092
$ws.Range("I136").Value = 58824980
$ws.Range("J136").Value = 8216.143
$ws.Range("K136").Value = 176474940
$ws.Range("L136").Value = 24648.429
$ws.Range("M136").Value = -176472390
$ws.Range("N136").Value = -29748.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1466.6666
$ws.Range("I4").Value = 1466.6666
$ws.Range("K4").Value = 1466.6666
$ws.Range("M4").Value = -1351.6666

$ws.Range("H20").Value = 3382.4546
$ws.Range("I20").Value = 4028.4285
$ws.Range("J20").Value = 2252
$ws.Range("K20").Value = 4028.4285
$ws.Range("L20").Value = 2252
$ws.Range("M20").Value = -3781.4285
$ws.Range("N20").Value = -2746

$ws.Range("H86").Value = 12843
$ws.Range("I86").Value = 6268.9165
$ws.Range("J86").Value = 44398.6
$ws.Range("K86").Value = 6268.9165
$ws.Range("L86").Value = 44398.6
$ws.Range("M86").Value = -5145.9165
$ws.Range("N86").Value = -46644.6

$ws.Range("H89").Value = 12843
$ws.Range("I89").Value = 6268.9165
$ws.Range("J89").Value = 44398.6
$ws.Range("K89").Value = 31344.5825
$ws.Range("L89").Value = 221993
$ws.Range("M89").Value = -25728.5825
$ws.Range("N89").Value = -233225

$ws.Range("H94").Value = 1652
$ws.Range("I94").Value = 1425.7693
$ws.Range("K94").Value = 1425.7693
$ws.Range("M94").Value = -974.7692999999999

$ws.Range("H99").Value = 3513
$ws.Range("I99").Value = 1816.3
$ws.Range("K99").Value = 1816.3
$ws.Range("M99").Value = -318.3

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H134").Value = 2437.8293
$ws.Range("I134").Value = 2329.9211
$ws.Range("K134").Value = 6989.763300000001
$ws.Range("M134").Value = -4454.763300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 449.6316
$ws.Range("I7").Value = 321.57144
$ws.Range("K7").Value = 321.57144
$ws.Range("M7").Value = -208.57144

$ws.Range("H141").Value = 376997.16
$ws.Range("I141").Value = 75098.336
$ws.Range("J141").Value = 467566.8
$ws.Range("K141").Value = 75098.336
$ws.Range("L141").Value = 467566.8
$ws.Range("M141").Value = -69918.336
$ws.Range("N141").Value = -477926.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1570
$ws.Range("J5").Value = 6000
$ws.Range("L5").Value = 18000
$ws.Range("N5").Value = -18224

$ws.Range("H114").Value = 1499
$ws.Range("I114").Value = 1499
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 4497
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H131").Value = 32050.406
$ws.Range("J131").Value = 5223.769
$ws.Range("L131").Value = 15671.307
$ws.Range("N131").Value = -25751.307

$ws.Range("H133").Value = 11267.454
$ws.Range("J133").Value = 19500.25
$ws.Range("L133").Value = 58500.75
$ws.Range("N133").Value = -68620.75

$ws.Range("H135").Value = 1570
$ws.Range("J135").Value = 6000
$ws.Range("L135").Value = 54000
$ws.Range("N135").Value = -59070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4504
$ws.Range("I70").Value = 4008
$ws.Range("K70").Value = 4008
$ws.Range("M70").Value = -3738

$ws.Range("H73").Value = 4504
$ws.Range("I73").Value = 4008
$ws.Range("K73").Value = 4008
$ws.Range("M73").Value = -3072

$ws.Range("H97").Value = 1563.0526
$ws.Range("I97").Value = 508.8
$ws.Range("J97").Value = 2734.4443
$ws.Range("K97").Value = 508.8
$ws.Range("L97").Value = 2734.4443
$ws.Range("M97").Value = -12.80000000000001
$ws.Range("N97").Value = -3726.4443

$ws.Range("H122").Value = 2299.75
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 2666.3333
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 7998.999899999999
$ws.Range("M122").Value = -1150
$ws.Range("N122").Value = -12898.9999

$ws.Range("H132").Value = 2895.0938
$ws.Range("I132").Value = 2793.9285
$ws.Range("K132").Value = 8381.7855
$ws.Range("M132").Value = -5851.7855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H47").Value = 34495
$ws.Range("J47").Value = 34495
$ws.Range("L47").Value = 34495

$ws.Range("H52").Value = 34495
$ws.Range("J52").Value = 34495
$ws.Range("L52").Value = 34495

$ws.Range("H107").Value = 10508.75
$ws.Range("I107").Value = 10508.75
$ws.Range("K107").Value = 10508.75
$ws.Range("M107").Value = -8588.75

$ws.Range("H132").Value = 233334860
$ws.Range("I132").Value = 1043
$ws.Range("K132").Value = 3129
$ws.Range("M132").Value = -599

$ws.Range("H136").Value = 1832.3096
$ws.Range("I136").Value = 1832.3096
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5496.9288
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4004.9636
$ws.Range("I132").Value = 4042.0637
$ws.Range("J132").Value = 3787
$ws.Range("K132").Value = 12126.1911
$ws.Range("L132").Value = 11361
$ws.Range("M132").Value = -9596.1911
$ws.Range("N132").Value = -16421
